$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-03 Monday" "2025-02-04 Tuesday"

Replace-Text "938×9=8442" "153×7=1071"
Replace-Text "744×3=2232" "437×3=1311"
Replace-Text "230×4=920" "766×4=3064"
Replace-Text "113×3=339" "255×7=1785"
Replace-Text "732×3=2196" "403×3=1209"
Replace-Text "678×4=2712" "425×5=2125"
Replace-Text "432×2=864" "903×4=3612"
Replace-Text "984×7=6888" "413×7=2891"
Replace-Text "288×2=576" "955×7=6685"
Replace-Text "783×8=6264" "135×6=810"
Replace-Text "739×4=2956" "879×7=6153"
Replace-Text "635×3=1905" "406×4=1624"
Replace-Text "281×4=1124" "101×2=202"
Replace-Text "646×7=4522" "131×2=262"
Replace-Text "825×4=3300" "296×9=2664"
Replace-Text "777×2=1554" "972×9=8748"
Replace-Text "235×2=470" "642×5=3210"
Replace-Text "629×3=1887" "921×7=6447"
Replace-Text "918×7=6426" "289×3=867"
Replace-Text "475×5=2375" "746×3=2238"
Replace-Text "376×9=3384" "682×5=3410"
Replace-Text "242×2=484" "498×5=2490"
Replace-Text "630×8=5040" "534×9=4806"
Replace-Text "757×8=6056" "179×3=537"
Replace-Text "297×2=594" "706×4=2824"
